$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Insert the new "2022-Q3" sheet by duplicating the existing "2022-Q2"
#    sheet (this preserves cell styles exactly) and placing the copy right
#    before it.
# ---------------------------------------------------------------------------
$srcQ2 = $wb.Worksheets.Item("2022-Q2")
$srcQ2.Copy($srcQ2)

# The newly created copy is now the sheet immediately before the original
# "2022-Q2" sheet (named "2022-Q2 (2)" by default). Rename it.
$newQ3 = $wb.Worksheets.Item("2022-Q2 (2)")
$newQ3.Name = "2022-Q3"

# Update the fund figures on the new sheet to the Q3 numbers (fund codes /
# names stay the same, only the holdings data changes). These columns are
# stored as text in the source workbook, so force a text number format
# before assigning so leading/format is preserved verbatim.
$newQ3.Range("D2:G3").NumberFormat = "@"

$newQ3.Range("D2").Value = "22.94"
$newQ3.Range("E2").Value = "77.45"
$newQ3.Range("F2").Value = "3.06"
$newQ3.Range("G2").Value = "0.7020"

$newQ3.Range("D3").Value = "4.07"
$newQ3.Range("E3").Value = "77.45"
$newQ3.Range("F3").Value = "3.06"
$newQ3.Range("G3").Value = "0.1245"

# ---------------------------------------------------------------------------
# 2. Update the "总计" (totals) summary sheet: insert a new row for 2022-Q3
#    at the top of the data and shift the rest down (the previously-last
#    "2020-Q4" row keeps its values, just moves to row 7).
# ---------------------------------------------------------------------------
$tot = $wb.Worksheets.Item("总计")

# A new row (7) is being added below the previously-last row (6); copy the
# cell format from A6 so the new A7 cell keeps the bordered/bold "index"
# column style used by the rest of column A.
$tot.Range("A6").Copy($tot.Range("A7"))

$rows = @(
    @{ A = 0; B = "2022-Q3"; C = 2;  D = 0.83 },
    @{ A = 1; B = "2022-Q2"; C = 2;  D = 0.55 },
    @{ A = 2; B = "2021-Q3"; C = 2;  D = 0.1 },
    @{ A = 3; B = "2021-Q2"; C = 14; D = 4.94 },
    @{ A = 4; B = "2021-Q1"; C = 7;  D = 3.52 },
    @{ A = 5; B = "2020-Q4"; C = 5;  D = 3.22 }
)

for ($i = 0; $i -lt $rows.Length; $i++) {
    $r = $i + 2
    $data = $rows[$i]
    $tot.Cells.Item($r, 1).Value = $data.A
    $tot.Cells.Item($r, 2).Value = $data.B
    $tot.Cells.Item($r, 3).Value = $data.C
    $tot.Cells.Item($r, 4).Value = $data.D
}
